# "Apresenta informação extra de ação" use-case sheet:
# add a second exception/alternative block ("Alternativa 2 ...") that
# mirrors the existing "Exceção 1 ..." block (rows 11-14), appended as
# rows 15-18, and wire it into the merged-cell / label layout the same
# way the existing block is wired.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reserve the merge for the new block's label column first so that when
# the formatting of the existing block is copied onto it, the merged
# region picks up the very same border/fill/alignment style already used
# by B11:B14 (merging *after* the copy causes Excel to recompute a brand
# new border combination for the merged area instead of reusing the
# existing one).
$ws.Range("B15:B18").Merge() | Out-Null

# Clone the look (fills, borders, alignment, number formats - not the
# values) of the existing "Exceção 1" block onto the four new rows.
$ws.Range("B11:D14").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New block's own text.
$ws.Range("B15").Value = " Alternativa 2 [requisita repetição de esclarecimento] (Passo 1)"
$ws.Range("D15").Value = "1.1 Regressa a 1"

# Mirror row 12's (empty C/D apart from the "1.2 Impede..." note) blank
# continuation rows - C16/C17/C18 and D17/D18 stay blank like their
# counterparts in the source block; D16 likewise stays blank (no third
# sub-step for the new alternative).

# Reflect the new bottom of the table / current selection.
$ws.Range("C22").Select() | Out-Null
